# Natmi following Dr Hou advice
# Rewrites the LR-pairs table (Icam1-Itgam) to the expanded 20-row grid
# (5 sender clusters x 4 target clusters), per the updated NATMI output.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$nRows = 20
$nCols = 20
$arr = New-Object 'object[,]' $nRows,$nCols

# Row 2: ECs -> FAPs
$arr[0,0] = "ECs"
$arr[0,1] = "Icam1"
$arr[0,2] = "Itgam"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 22.481209
$arr[0,7] = 67.443627
$arr[0,8] = 0.1656226259370683
$arr[0,9] = 0.166106832923046
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 0.9636300000000001
$arr[0,13] = 2.89089
$arr[0,14] = 0.00463428088190967
$arr[0,15] = 0.004636691828827675
$arr[0,16] = 21.66356742867001
$arr[0,17] = 194.97210685803
$arr[0,18] = 0.0007675417689918324
$arr[0,19] = 0.0007701861949267311

# Row 3: ECs -> M1
$arr[1,0] = "ECs"
$arr[1,1] = "Icam1"
$arr[1,2] = "Itgam"
$arr[1,3] = "M1"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 22.481209
$arr[1,7] = 67.443627
$arr[1,8] = 0.1656226259370683
$arr[1,9] = 0.166106832923046
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 71.625121
$arr[1,13] = 214.875363
$arr[1,14] = 0.3444588990740914
$arr[1,15] = 0.344638101013349
$arr[1,16] = 1610.219314851289
$arr[1,17] = 14491.9738336616
$arr[1,18] = 0.05705018739204262
$arr[1,19] = 0.05724674346394021

# Row 4: ECs -> M2
$arr[2,0] = "ECs"
$arr[2,1] = "Icam1"
$arr[2,2] = "Itgam"
$arr[2,3] = "M2"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 22.481209
$arr[2,7] = 67.443627
$arr[2,8] = 0.1656226259370683
$arr[2,9] = 0.166106832923046
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 135.0220616666666
$arr[2,13] = 405.066185
$arr[2,14] = 0.6493469059886695
$arr[2,15] = 0.6496847234325412
$arr[2,16] = 3035.459187939221
$arr[2,17] = 27319.132691453
$arr[2,18] = 0.1075465397139541
$arr[2,19] = 0.1079170718078645

# Row 5: ECs -> sCs
$arr[3,0] = "ECs"
$arr[3,1] = "Icam1"
$arr[3,2] = "Itgam"
$arr[3,3] = "sCs"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 22.481209
$arr[3,7] = 67.443627
$arr[3,8] = 0.1656226259370683
$arr[3,9] = 0.166106832923046
$arr[3,10] = 1
$arr[3,11] = 0.5
$arr[3,12] = 0.324361
$arr[3,13] = 0.648722
$arr[3,14] = 0.001559914055329434
$arr[3,15] = 0.001040483725282092
$arr[3,16] = 7.292027432449001
$arr[3,17] = 43.75216459469401
$arr[3,18] = 0.0002583570620798021
$arr[3,19] = 0.0001728314563145809

# Row 6: FAPs -> FAPs
$arr[4,0] = "FAPs"
$arr[4,1] = "Icam1"
$arr[4,2] = "Itgam"
$arr[4,3] = "FAPs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 29.50180766666667
$arr[4,7] = 88.50542300000001
$arr[4,8] = 0.2173444878184117
$arr[4,9] = 0.2179799065528387
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 0.9636300000000001
$arr[4,13] = 2.89089
$arr[4,14] = 0.00463428088190967
$arr[4,15] = 0.004636691828827675
$arr[4,16] = 28.42882692183001
$arr[4,17] = 255.85944229647
$arr[4,18] = 0.001007235404685314
$arr[4,19] = 0.001010705651562167

# Row 7: FAPs -> M1
$arr[5,0] = "FAPs"
$arr[5,1] = "Icam1"
$arr[5,2] = "Itgam"
$arr[5,3] = "M1"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 29.50180766666667
$arr[5,7] = 88.50542300000001
$arr[5,8] = 0.2173444878184117
$arr[5,9] = 0.2179799065528387
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 71.625121
$arr[5,13] = 214.875363
$arr[5,14] = 0.3444588990740914
$arr[5,15] = 0.344638101013349
$arr[5,16] = 2113.070543843728
$arr[5,17] = 19017.63489459355
$arr[5,18] = 0.07486624299375237
$arr[5,19] = 0.07512418105343761

# Row 8: FAPs -> M2
$arr[6,0] = "FAPs"
$arr[6,1] = "Icam1"
$arr[6,2] = "Itgam"
$arr[6,3] = "M2"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 29.50180766666667
$arr[6,7] = 88.50542300000001
$arr[6,8] = 0.2173444878184117
$arr[6,9] = 0.2179799065528387
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 135.0220616666666
$arr[6,13] = 405.066185
$arr[6,14] = 0.6493469059886695
$arr[6,15] = 0.6496847234325412
$arr[6,16] = 3983.394894046806
$arr[6,17] = 35850.55404642125
$arr[6,18] = 0.1411319706985777
$arr[6,19] = 0.1416182153026322

# Row 9: FAPs -> sCs
$arr[7,0] = "FAPs"
$arr[7,1] = "Icam1"
$arr[7,2] = "Itgam"
$arr[7,3] = "sCs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 29.50180766666667
$arr[7,7] = 88.50542300000001
$arr[7,8] = 0.2173444878184117
$arr[7,9] = 0.2179799065528387
$arr[7,10] = 1
$arr[7,11] = 0.5
$arr[7,12] = 0.324361
$arr[7,13] = 0.648722
$arr[7,14] = 0.001559914055329434
$arr[7,15] = 0.001040483725282092
$arr[7,16] = 9.569235836567668
$arr[7,17] = 57.41541501940601
$arr[7,18] = 0.0003390387213963173
$arr[7,19] = 0.0002268045452067399

# Row 10: M1 -> FAPs
$arr[8,0] = "M1"
$arr[8,1] = "Icam1"
$arr[8,2] = "Itgam"
$arr[8,3] = "FAPs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 42.765269
$arr[8,7] = 128.295807
$arr[8,8] = 0.3150585073376215
$arr[8,9] = 0.3159795984589671
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 0.9636300000000001
$arr[8,13] = 2.89089
$arr[8,14] = 0.00463428088190967
$arr[8,15] = 0.004636691828827675
$arr[8,16] = 41.20989616647
$arr[8,17] = 370.88906549823
$arr[8,18] = 0.001460069617237737
$arr[8,19] = 0.001465100022250942

# Row 11: M1 -> M1
$arr[9,0] = "M1"
$arr[9,1] = "Icam1"
$arr[9,2] = "Itgam"
$arr[9,3] = "M1"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 42.765269
$arr[9,7] = 128.295807
$arr[9,8] = 0.3150585073376215
$arr[9,9] = 0.3159795984589671
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 71.625121
$arr[9,13] = 214.875363
$arr[9,14] = 0.3444588990740914
$arr[9,15] = 0.344638101013349
$arr[9,16] = 3063.067566722548
$arr[9,17] = 27567.60810050294
$arr[9,18] = 0.1085247065814437
$arr[9,19] = 0.1088986087718589

# Row 12: M1 -> M2
$arr[10,0] = "M1"
$arr[10,1] = "Icam1"
$arr[10,2] = "Itgam"
$arr[10,3] = "M2"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 42.765269
$arr[10,7] = 128.295807
$arr[10,8] = 0.3150585073376215
$arr[10,9] = 0.3159795984589671
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 135.0220616666666
$arr[10,13] = 405.066185
$arr[10,14] = 0.6493469059886695
$arr[10,15] = 0.6496847234325412
$arr[10,16] = 5774.254788109587
$arr[10,17] = 51968.29309298629
$arr[10,18] = 0.204582266945093
$arr[10,19] = 0.2052871180351394

# Row 13: M1 -> sCs
$arr[11,0] = "M1"
$arr[11,1] = "Icam1"
$arr[11,2] = "Itgam"
$arr[11,3] = "sCs"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 42.765269
$arr[11,7] = 128.295807
$arr[11,8] = 0.3150585073376215
$arr[11,9] = 0.3159795984589671
$arr[11,10] = 1
$arr[11,11] = 0.5
$arr[11,12] = 0.324361
$arr[11,13] = 0.648722
$arr[11,14] = 0.001559914055329434
$arr[11,15] = 0.001040483725282092
$arr[11,16] = 13.871385418109
$arr[11,17] = 83.228312508654
$arr[11,18] = 0.0004914641938470672
$arr[11,19] = 0.0003287716297177256

# Row 14: M2 -> FAPs
$arr[12,0] = "M2"
$arr[12,1] = "Icam1"
$arr[12,2] = "Itgam"
$arr[12,3] = "FAPs"
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 39.80222300000001
$arr[12,7] = 119.406669
$arr[12,8] = 0.2932292783449848
$arr[12,9] = 0.2940865504976542
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 0.9636300000000001
$arr[12,13] = 2.89089
$arr[12,14] = 0.00463428088190967
$arr[12,15] = 0.004636691828827675
$arr[12,16] = 38.35461614949001
$arr[12,17] = 345.19154534541
$arr[12,18] = 0.001358906838650332
$arr[12,19] = 0.001363588705660591

# Row 15: M2 -> M1
$arr[13,0] = "M2"
$arr[13,1] = "Icam1"
$arr[13,2] = "Itgam"
$arr[13,3] = "M1"
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 39.80222300000001
$arr[13,7] = 119.406669
$arr[13,8] = 0.2932292783449848
$arr[13,9] = 0.2940865504976542
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 71.625121
$arr[13,13] = 214.875363
$arr[13,14] = 0.3444588990740914
$arr[13,15] = 0.344638101013349
$arr[13,16] = 2850.839038443983
$arr[13,17] = 25657.55134599585
$arr[13,18] = 0.1010054343950038
$arr[13,19] = 0.1013534302970779

# Row 16: M2 -> M2
$arr[14,0] = "M2"
$arr[14,1] = "Icam1"
$arr[14,2] = "Itgam"
$arr[14,3] = "M2"
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 39.80222300000001
$arr[14,7] = 119.406669
$arr[14,8] = 0.2932292783449848
$arr[14,9] = 0.2940865504976542
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 135.0220616666666
$arr[14,13] = 405.066185
$arr[14,14] = 0.6493469059886695
$arr[14,15] = 0.6496847234325412
$arr[14,16] = 5374.178208376418
$arr[14,17] = 48367.60387538777
$arr[14,18] = 0.1904075246386062
$arr[14,19] = 0.1910635392252985

# Row 17: M2 -> sCs
$arr[15,0] = "M2"
$arr[15,1] = "Icam1"
$arr[15,2] = "Itgam"
$arr[15,3] = "sCs"
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 39.80222300000001
$arr[15,7] = 119.406669
$arr[15,8] = 0.2932292783449848
$arr[15,9] = 0.2940865504976542
$arr[15,10] = 1
$arr[15,11] = 0.5
$arr[15,12] = 0.324361
$arr[15,13] = 0.648722
$arr[15,14] = 0.001559914055329434
$arr[15,15] = 0.001040483725282092
$arr[15,16] = 12.910288854503
$arr[15,17] = 77.46173312701801
$arr[15,18] = 0.0004574124727244485
$arr[15,19] = 0.0003059922696171593

# Row 18: sCs -> FAPs
$arr[16,0] = "sCs"
$arr[16,1] = "Icam1"
$arr[16,2] = "Itgam"
$arr[16,3] = "FAPs"
$arr[16,4] = 2
$arr[16,5] = 1
$arr[16,6] = 1.1870385
$arr[16,7] = 2.374077
$arr[16,8] = 0.00874510056191367
$arr[16,9] = 0.005847111567493934
$arr[16,10] = 3
$arr[16,11] = 1
$arr[16,12] = 0.9636300000000001
$arr[16,13] = 2.89089
$arr[16,14] = 0.00463428088190967
$arr[16,15] = 0.004636691828827675
$arr[16,16] = 1.143865909755
$arr[16,17] = 6.863195458530001
$arr[16,18] = 0.00004052725234445403
$arr[16,19] = 0.0000271112544272429

# Row 19: sCs -> M1
$arr[17,0] = "sCs"
$arr[17,1] = "Icam1"
$arr[17,2] = "Itgam"
$arr[17,3] = "M1"
$arr[17,4] = 2
$arr[17,5] = 1
$arr[17,6] = 1.1870385
$arr[17,7] = 2.374077
$arr[17,8] = 0.00874510056191367
$arr[17,9] = 0.005847111567493934
$arr[17,10] = 3
$arr[17,11] = 1
$arr[17,12] = 71.625121
$arr[17,13] = 214.875363
$arr[17,14] = 0.3444588990740914
$arr[17,15] = 0.344638101013349
$arr[17,16] = 85.0217761941585
$arr[17,17] = 510.130657164951
$arr[17,18] = 0.003012327711849002
$arr[17,19] = 0.002015137427034296

# Row 20: sCs -> M2
$arr[18,0] = "sCs"
$arr[18,1] = "Icam1"
$arr[18,2] = "Itgam"
$arr[18,3] = "M2"
$arr[18,4] = 2
$arr[18,5] = 1
$arr[18,6] = 1.1870385
$arr[18,7] = 2.374077
$arr[18,8] = 0.00874510056191367
$arr[18,9] = 0.005847111567493934
$arr[18,10] = 3
$arr[18,11] = 1
$arr[18,12] = 135.0220616666666
$arr[18,13] = 405.066185
$arr[18,14] = 0.6493469059886695
$arr[18,15] = 0.6496847234325412
$arr[18,16] = 160.2763855477075
$arr[18,17] = 961.658313286245
$arr[18,18] = 0.005678603992438417
$arr[18,19] = 0.003798779061606509

# Row 21: sCs -> sCs
$arr[19,0] = "sCs"
$arr[19,1] = "Icam1"
$arr[19,2] = "Itgam"
$arr[19,3] = "sCs"
$arr[19,4] = 2
$arr[19,5] = 1
$arr[19,6] = 1.1870385
$arr[19,7] = 2.374077
$arr[19,8] = 0.00874510056191367
$arr[19,9] = 0.005847111567493934
$arr[19,10] = 1
$arr[19,11] = 0.5
$arr[19,12] = 0.324361
$arr[19,13] = 0.648722
$arr[19,14] = 0.001559914055329434
$arr[19,15] = 0.001040483725282092
$arr[19,16] = 0.3850289948985001
$arr[19,17] = 1.540115979594
$arr[19,18] = 0.00001364160528179846
$arr[19,19] = 0.000006083824425886101

$ws.Range("A2:T21").Value2 = $arr
Write-Output "Wrote $nRows x $nCols block into A2:T21"
